# Add season-record columns (Wins / Losses / Ties) to the player table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new headers in AD1:AF1 ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the existing header formatting (bold / border / centered) by
# copying the format from an existing header cell (A1) onto the new ones,
# same as the rest of row 1.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# --- Data rows (2-47): same season record for every player row ---
$lastRow = 47
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 78   # AD = Wins
    $ws.Cells.Item($r, 31).Value = 83   # AE = Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF = Ties
}
